$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: "Best Time to Buy and Sell Stock" (LeetCode #121)
# Copy row 22's formatting down into the new row 23 first so the new
# row inherits the same fill/border/alignment/number-format styling.
$ws.Range("A22:H22").Copy()
$ws.Range("A23:H23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A23").Value = 121
$ws.Range("B23").Value = "Easy"
$ws.Range("C23").Value = "Best Time to Buy and Sell Stock"

# Add the hyperlink, then restore D22's exact cell formatting (Hyperlinks.Add
# tends to stamp its own style) so D23 matches the sheet's existing look.
$ws.Hyperlinks.Add($ws.Range("D23"), "https://shorturl.at/Ugo4p")
$ws.Range("D22").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F23").Value = "O(n)"
$ws.Range("G23").Value = "Non-intuitive logic, but forming the test cases and then testing the solution works"
$ws.Range("H23").Value = 45496

$ws.Range("C25").Select()
